# Updates the cryptos list (Price / Volume(1h) columns, plus a handful of
# Coin/Link row swaps caused by re-ranking) to match the latest scrape.
# Commit: "Updated cryptos list on Wed May  1 07:33:23 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D="58.095.85"; E="  -8.33%  "},
    @{Row=3; D="2.901.30"; E="  -8.46%  "},
    @{Row=4; E="  +0.16%  "},
    @{Row=5; D="548.07"; E="  -9.03%  "},
    @{Row=6; D="120.44"; E="  -11.35%  "},
    @{Row=7; E="  +0.33%  "},
    @{Row=8; D="2.897.29"; E="  -8.52%  "},
    @{Row=9; D="0.490"; E="  -4.41%  "},
    @{Row=10; D="0.125"; E="  -12.26%  "},
    @{Row=12; D="0.429"; E="  -5.44%  "},
    @{Row=13; E="  -11.62%  "},
    @{Row=14; D="31.30"; E="  -10.20%  "},
    @{Row=15; E="  -1.25%  "},
    @{Row=16; D="3.380.79"; E="  -8.44%  "},
    @{Row=17; B="WrappedBTC"; C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D="58.409.48"; E="  -7.81%  "},
    @{Row=18; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="2.892.31"; E="  -8.85%  "},
    @{Row=19; D="6.37"; E="  -3.28%  "},
    @{Row=20; D="414.00"; E="  -10.25%  "},
    @{Row=21; D="12.77"; E="  -8.58%  "},
    @{Row=22; D="0.648"; E="  -7.08%  "},
    @{Row=23; D="6.78"; E="  -11.41%  "},
    @{Row=24; D="12.46"; E="  -6.61%  "},
    @{Row=25; D="76.81"; E="  -7.71%  "},
    @{Row=26; E="  +0.31%  "},
    @{Row=27; E="  +0.16%  "},
    @{Row=28; D="2.44"; E="  -9.45%  "},
    @{Row=29; D="1.89"; E="  -9.27%  "},
    @{Row=30; D="6.96"; E="  -9.60%  "},
    @{Row=31; B="NEARProtocol"; C="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D="5.92"; E="  -12.45%  "},
    @{Row=32; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="24.47"; E="  -9.87%  "},
    @{Row=33; D="0.0931"; E="  -7.16%  "},
    @{Row=34; B="OKB"; C="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D="48.92"; E="  -4.51%  "},
    @{Row=35; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="5.35"; E="  -9.46%  "},
    @{Row=36; B="Mantle"; C="https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; D="0.892"; E="  -12.62%  "},
    @{Row=37; E="  -19.21%  "},
    @{Row=38; D="8.26"; E="  +1.50%  "},
    @{Row=39; D="0.0₃0612"; E="  -16.49%  "},
    @{Row=40; D="0.0342"; E="  -12.40%  "},
    @{Row=41; D="0.104"; E="  -7.52%  "},
    @{Row=42; D="2.606.81"; E="  -6.84%  "},
    @{Row=43; B="Bittensor"; C="https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; D="350.50"; E="  -10.37%  "},
    @{Row=44; B="dogwifhat"; C="https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; D="2.34"; E="  -10.96%  "},
    @{Row=45; E="  +0.00%  "},
    @{Row=46; D="118.20"; E="  -6.04%  "},
    @{Row=47; D="0.226"; E="  -10.03%  "},
    @{Row=48; D="0.105"; E="  -5.64%  "},
    @{Row=49; D="1.91"; E="  -9.99%  "},
    @{Row=50; D="22.49"; E="  -10.38%  "},
    @{Row=51; D="1.93"; E="  -11.20%  "},
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.ContainsKey("B")) { $ws.Range("B$row").Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C$row").Value = $u.C }

    if ($u.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $val = $u.D
        # The price column stores plain-looking numbers (e.g. "548.07") as
        # TEXT, not numbers (other rows use "."-grouped thousands like
        # "58.095.85", which Excel would never auto-parse as a number).
        # Force text formatting first so the literal string is preserved,
        # then restore the default "Normal" style so no stray formatting
        # is left behind on the cell.
        $looksNumeric = $val -match '^-?[0-9]+(\.[0-9]+)?$'
        if ($looksNumeric) {
            $cell.NumberFormat = "@"
            $cell.Value = $val
            $cell.Style = "Normal"
        } else {
            $cell.Value = $val
        }
    }

    if ($u.ContainsKey("E")) { $ws.Range("E$row").Value = $u.E }
}

Write-Output "Applied $($updates.Count) row updates"
